# Update "Chức năng nhiều đội phó quản lý 1 cán bộ"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Trim the two merged "Đơn vị/Vị trí công tác" blocks down to their new
#    sizes by deleting rows *inside* each merged block (so Excel shrinks the
#    existing merge in place instead of us having to unmerge/re-merge and
#    picking up stray border styling):
#      - "Nghiệp vụ" block was D5:D10 (6 rows) -> keep 4 rows -> delete 2
#        rows from inside it (old rows 9:10).
#      - "Văn phòng" block was D11:D14 (4 rows) -> keep 2 rows -> delete 2
#        rows from inside it (old rows 11:12, i.e. the first two rows of
#        that block, now sitting right after the shrunk first block).
$ws.Rows("9:10").Delete()
$ws.Rows("11:12").Delete()

# 2) Drop the "Số ngày nghỉ không phép" column (old column G); everything to
#    the right shifts left by one (M -> L).
$ws.Columns("G:G").Delete()

# 3) Title + subtitle
$ws.Range("A1").Value = "TỔNG HỢP Kết quả đánh giá, xếp loại chất lượng công chức"
$ws.Range("A2").Value = "Tháng 03/2025"

# 4) Header row (row 4)
$ws.Range("A4").Value = "STT"
$ws.Range("B4").Value = "Họ và tên"
$ws.Range("C4").Value = "Chức vụ"
$ws.Range("D4").Value = "Đơn vị/Vị trí công tác"
$ws.Range("E4").Value = "Số ngày làm việc thực tế"
$ws.Range("F4").Value = "Số ngày nghỉ có phép"
$ws.Range("G4").Value = "Số lần vi phạm quy chế, quy định"
$ws.Range("H4").Value = "Hình thức kỷ luật"
$ws.Range("I4").Value = "Tự xếp loại"
$ws.Range("J4").Value = "% mức độ hoàn thành nhiệm vụ"
$ws.Range("K4").Value = "Mức xếp loại của Lãnh đạo"
$ws.Range("L4").Value = "Tổng Nhiệm Vụ"

# 5) Data rows 5-10
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Hoàng Văn Nam"
$ws.Range("C5").Value = "Đội trưởng"
$ws.Range("D5").Value = "Nghiệp vụ"
$ws.Range("I5").Value = "B"
$ws.Range("J5").Value = 100.0
$ws.Range("K5").Value = "B"
$ws.Range("L5").Value = 1

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Đội phó 2"
$ws.Range("C6").Value = "Đội phó"
$ws.Range("E6").Value = 20
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = "D"
$ws.Range("K6").Value = "D"
$ws.Range("L6").Value = 1

$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "Đội Phó"
$ws.Range("C7").Value = "Đội phó"
$ws.Range("I7").Value = "C"
$ws.Range("K7").Value = "B"

$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "Nguyễn Văn Hải"
$ws.Range("C8").Value = "Đội phó"
$ws.Range("E8").Value = ""
$ws.Range("H8").Value = ""
$ws.Range("I8").Value = "C"
$ws.Range("K8").Value = "B"

$ws.Range("A9").Value = 5
$ws.Range("B9").Value = "Phan Nghĩa"
$ws.Range("C9").Value = "Đội phó"
$ws.Range("D9").Value = "Văn phòng"
$ws.Range("I9").Value = "C"
$ws.Range("K9").Value = "C"

$ws.Range("A10").Value = 6
$ws.Range("B10").Value = "Nguyễn Hiếu"
$ws.Range("C10").Value = "Công chức"
$ws.Range("E10").Value = 23
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = "B"
$ws.Range("J10").Value = 100.0
$ws.Range("K10").Value = "B"
$ws.Range("L10").Value = 1

# 6) Refresh the sheet's remembered selection to match the new, smaller
#    used range (A4:M14 -> A4:L10).
$ws.Range("A4:L10").Select() | Out-Null
